# Updated cryptos list (Price / Volume(1h) columns) for Sheet1 rows 2-51.
# Values that look purely numeric (e.g. "683.69", "0.999") are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original inline-string cells, instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.273.75'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '3.681.45'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'683.69"
$ws.Range('E5').Value = '  -3.11%  '
$ws.Range('D6').Value = "'162.10"
$ws.Range('E6').Value = '  -5.62%  '
$ws.Range('D7').Value = '3.680.77'
$ws.Range('E7').Value = '  -3.53%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.500"
$ws.Range('E9').Value = '  -4.85%  '
$ws.Range('E10').Value = '  -8.46%  '
$ws.Range('D11').Value = "'7.39"
$ws.Range('E11').Value = '  -3.29%  '
$ws.Range('D12').Value = "'0.445"
$ws.Range('E12').Value = '  -4.05%  '
$ws.Range('E13').Value = '  -5.26%  '
$ws.Range('D14').Value = "'33.66"
$ws.Range('E14').Value = '  -6.57%  '
$ws.Range('D15').Value = '4.301.85'
$ws.Range('E15').Value = '  -3.52%  '
$ws.Range('D16').Value = '3.682.06'
$ws.Range('E16').Value = '  -3.54%  '
$ws.Range('D17').Value = '69.332.86'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = "'16.35"
$ws.Range('E19').Value = '  -6.62%  '
$ws.Range('D20').Value = "'6.62"
$ws.Range('E20').Value = '  -7.78%  '
$ws.Range('D21').Value = "'481.48"
$ws.Range('E21').Value = '  -7.33%  '
$ws.Range('D22').Value = "'9.94"
$ws.Range('E22').Value = '  -6.81%  '
$ws.Range('D23').Value = "'0.666"
$ws.Range('E23').Value = '  -8.13%  '
$ws.Range('D24').Value = "'80.24"
$ws.Range('E24').Value = '  -5.23%  '
$ws.Range('D25').Value = '3.825.88'
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('E26').Value = '  -9.98%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = "'11.44"
$ws.Range('E28').Value = '  -5.19%  '
$ws.Range('D29').Value = "'9.50"
$ws.Range('E29').Value = '  -8.84%  '
$ws.Range('E30').Value = '  -10.33%  '
$ws.Range('D31').Value = "'2.71"
$ws.Range('E31').Value = '  -10.52%  '
$ws.Range('D32').Value = "'2.08"
$ws.Range('E32').Value = '  -7.65%  '
$ws.Range('D33').Value = "'6.84"
$ws.Range('E33').Value = '  -7.50%  '
$ws.Range('D34').Value = "'0.168"
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('D35').Value = "'27.12"
$ws.Range('E35').Value = '  -7.25%  '
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Value = '3.650.65'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('E38').Value = '  -7.67%  '
$ws.Range('D39').Value = "'6.30"
$ws.Range('E39').Value = '  +5.86%  '
$ws.Range('D40').Value = "'2.34"
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = "'0.0935"
$ws.Range('E41').Value = '  -7.93%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = "'0.948"
$ws.Range('E44').Value = '  -7.28%  '
$ws.Range('D45').Value = "'162.63"
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('D46').Value = "'48.33"
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').Value = "'2.83"
$ws.Range('E47').Value = '  -13.47%  '
$ws.Range('D48').Value = "'29.96"
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('D49').Value = "'0.000288"
$ws.Range('E49').Value = '  -7.88%  '
$ws.Range('E50').Value = '  -1.44%  '
$ws.Range('D51').Value = "'1.12"
$ws.Range('E51').Value = '  -3.28%  '

# Clear the quote-prefix style hint picked up above so cell styling
# matches the original (text stays text; no stray number format).
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
